# daily auto push: 2026-01-20 19:08 UTC
# Insert two new rows into the daily log at row 665 (pushing the existing
# 2026/12/29.. rows down by two) and fill them with the 2026/01/20 /
# 2026/01/21 entries that were missing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows starting at row 665; everything from the old
# row 665 onward shifts down to row 667 onward.
$ws.Range("A665:A666").EntireRow.Insert()

# Row 665: 2026/01/20, 火, 22, 19
# Force text format before writing so the date-like string is not
# auto-converted into a real date serial, then restore the default
# "Normal" style so no extra formatting is left behind on the cell.
$ws.Range("A665").NumberFormat = "@"
$ws.Range("A665").Value = "2026/01/20"
$ws.Range("A665").Style = "Normal"
$ws.Range("B665").Value = "火"
$ws.Range("C665").Value = 22
$ws.Range("D665").Value = 19

# Row 666: 2026/01/21, 水, 2, 21
$ws.Range("A666").NumberFormat = "@"
$ws.Range("A666").Value = "2026/01/21"
$ws.Range("A666").Style = "Normal"
$ws.Range("B666").Value = "水"
$ws.Range("C666").Value = 2
$ws.Range("D666").Value = 21
